# Fill in the evidence workbook with the real on-chain data that replaces
# the placeholder "hint" text in the A1 and A2 submission sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "A1": TxHash / ClassID ------------------------------------------------
$wsA1 = $wb.Worksheets.Item("A1")

# Carry over the header cells' formatting (Arial + thin border) onto the
# answer row, then overwrite the placeholder hint text with the real values.
$wsA1.Range("A1").Copy()
$wsA1.Range("A2").PasteSpecial(-4122)
$wsA1.Range("B1").Copy()
$wsA1.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsA1.Range("A2").Value = "B37ECC9D29AFD159DA1E3A7120F44A53F8539A05ECE567AAD4B6E2DBD5E3A0E1"
$wsA1.Range("B2").Value = "jav1xgon1"

$wsA1.Range("B16").Select()

# --- Sheet "A2": TxHash / ClassID / NFTID ----------------------------------------
$wsA2 = $wb.Worksheets.Item("A2")

$wsA2.Range("A1").Copy()
$wsA2.Range("A2").PasteSpecial(-4122)
$wsA2.Range("B1").Copy()
$wsA2.Range("B2").PasteSpecial(-4122)
$wsA2.Range("C1").Copy()
$wsA2.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsA2.Range("A2").Value = "43DA6B57F90E244DBAC3E0833EF58906BC8279B2E513BF2EBE36320943FB867E"
$wsA2.Range("B2").Value = "jav1xgon1"
$wsA2.Range("C2").Value = "jav1xgon"

$wsA2.Range("E11").Select()

# --- Restore the original active sheet/selection --------------------------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Select()
$wsInfo.Range("F2").Select()
